# Weekly crypto snapshot refresh (GitHub Actions data pull).
# Updates Price (D) / Volume(1h) (E) figures for most rows, and for two
# pairs of coins (rows 30/31, 41/42, 48/49) the ranking swapped places
# so the coin name / link / price / volume moved to the other row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as plain TEXT (e.g. '64.212.85', '1.00', '0.999')
# so decimal points / trailing zeros are preserved exactly as scraped.
# Excel's COM layer auto-converts numeric-looking strings assigned via
# .Value into real numbers, so mark the range as Text first, then strip
# the formatting back off afterwards so the cells keep their original
# (default) style while remaining Text-typed.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '64.212.85'
$ws.Range("E2").Value = '  +1.71%  '
# Row 3
$ws.Range("D3").Value = '2.627.32'
$ws.Range("E3").Value = '  -0.02%  '
# Row 5
$ws.Range("D5").Value = '602.66'
$ws.Range("E5").Value = '  +0.09%  '
# Row 6
$ws.Range("D6").Value = '150.91'
# Row 7
$ws.Range("E7").Value = '  +0.04%  '
# Row 8
$ws.Range("E8").Value = '  +0.60%  '
# Row 9
$ws.Range("E9").Value = '  +2.32%  '
# Row 10
$ws.Range("D10").Value = '5.73'
$ws.Range("E10").Value = '  +2.48%  '
# Row 11
$ws.Range("D11").Value = '0.386'
$ws.Range("E11").Value = '  +6.51%  '
# Row 12
$ws.Range("E12").Value = '  -0.76%  '
# Row 13
$ws.Range("D13").Value = '27.67'
$ws.Range("E13").Value = '  +2.03%  '
# Row 14
$ws.Range("D14").Value = '3.099.75'
$ws.Range("E14").Value = '  +0.02%  '
# Row 15
$ws.Range("D15").Value = '64.065.38'
# Row 16
$ws.Range("D16").Value = '0.0000149'
$ws.Range("E16").Value = '  +3.72%  '
# Row 17
$ws.Range("D17").Value = '2.625.03'
$ws.Range("E17").Value = '  +0.37%  '
# Row 18
$ws.Range("D18").Value = '12.20'
$ws.Range("E18").Value = '  +8.16%  '
# Row 19
$ws.Range("D19").Value = '4.66'
$ws.Range("E19").Value = '  +3.93%  '
# Row 20
$ws.Range("D20").Value = '352.12'
$ws.Range("E20").Value = '  +3.50%  '
# Row 21
$ws.Range("E21").Value = '  +1.37%  '
# Row 22
$ws.Range("E22").Value = '  -0.15%  '
# Row 23
$ws.Range("E23").Value = '  +2.96%  '
# Row 24
$ws.Range("D24").Value = '66.68'
$ws.Range("E24").Value = '  +0.40%  '
# Row 25
$ws.Range("D25").Value = '1.76'
$ws.Range("E25").Value = '  +16.59%  '
# Row 26
$ws.Range("E26").Value = '  +4.98%  '
# Row 27
$ws.Range("D27").Value = '9.25'
$ws.Range("E27").Value = '  +6.97%  '
# Row 28
$ws.Range("D28").Value = '0.165'
$ws.Range("E28").Value = '  +1.74%  '
# Row 29
$ws.Range("D29").Value = '8.11'
$ws.Range("E29").Value = '  +3.50%  '
# Row 30
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.16%  '
# Row 31
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '537.22'
$ws.Range("E31").Value = '  -1.32%  '
# Row 32
$ws.Range("E32").Value = '  +10.08%  '
# Row 33
$ws.Range("D33").Value = '0.0₃0855'
$ws.Range("E33").Value = '  +6.61%  '
# Row 34
$ws.Range("E34").Value = '  +0.56%  '
# Row 35
$ws.Range("D35").Value = '5.29'
$ws.Range("E35").Value = '  -0.55%  '
# Row 36
$ws.Range("D36").Value = '167.82'
$ws.Range("E36").Value = '  +0.52%  '
# Row 37
$ws.Range("E37").Value = '  +7.79%  '
# Row 38
$ws.Range("E38").Value = '  +1.96%  '
# Row 39
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  +0.01%  '
# Row 40
$ws.Range("D40").Value = '19.55'
# Row 41
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.03%  '
# Row 42
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '170.04'
$ws.Range("E42").Value = '  +1.34%  '
# Row 43
$ws.Range("D43").Value = '39.95'
$ws.Range("E43").Value = '  +0.75%  '
# Row 44
$ws.Range("E44").Value = '  +5.57%  '
# Row 45
$ws.Range("E45").Value = '  +4.36%  '
# Row 46
$ws.Range("D46").Value = '21.56'
$ws.Range("E46").Value = '  -3.46%  '
# Row 47
$ws.Range("E47").Value = '  +1.34%  '
# Row 48
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '0.0247'
$ws.Range("E48").Value = '  +1.73%  '
# Row 49
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").Value = '2.00'
$ws.Range("E49").Value = '  +13.99%  '
# Row 50
$ws.Range("D50").Value = '0.0968'
$ws.Range("E50").Value = '  +1.22%  '
# Row 51
$ws.Range("D51").Value = '19.32'
$ws.Range("E51").Value = '  +4.17%  '

# Restore the default (unstyled) formatting now that the text values are
# safely stored, so only the cell VALUES differ from the original file.
$ws.Range("D2:D51").ClearFormats()
